# Applies the IFRS re-statement of 한국콜마홀딩스 financial figures:
#  - Rows 2-6 (actual years 2014-2018): replace D..AJ with corrected figures
#  - Rows 7-9 (estimate years 2019E-2021E): clear all D..AJ figures (unavailable)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{
        D=2225; E=277; F=327; G=341; H=274; I=195; J=79; K=3782; L=1487; M=2295;
        N=1988; O=307; P=82; Q=292; R=-818; S=564; T=188; U=104; V=973;
        W=12.46; X=12.34; Y=10.33; Z=8.33; AA=64.79000000000001; AB=2335.38;
        AC=1190; AD=36.29; AE=12204; AF=3.54; AG=100; AH=0.23; AI=8.35; AJ=16391930
    }
    3 = @{
        D=2996; E=451; F=542; G=109; H=10; I=-7; J=18; K=4995; L=2317; M=2678;
        N=2151; O=528; P=84; Q=265; R=-857; S=691; T=84; U=181; V=1588;
        W=15.05; X=0.35; Y=-0.36; Z=0.24; AA=86.51000000000001; AB=2477.56;
        AC=-45; AD=-1435.92; AE=12915; AF=5.03; AG=125; AH=0.19; AI=-280.38; AJ=16748936
    }
    4 = @{
        D=3304; E=487; F=612; G=621; H=510; I=366; J=144; K=5676; L=2485; M=3191;
        N=2528; O=663; P=84; Q=253; R=-407; S=222; T=131; U=122; V=1946;
        W=14.75; X=15.43; Y=15.67; Z=9.56; AA=77.88; AB=2920.84;
        AC=2188; AD=15.4; AE=15180; AF=2.22; AG=155; AH=0.46; AI=7.04; AJ=16748936
    }
    5 = @{
        D=4083; E=618; F=818; G=676; H=517; I=280; J=236; K=7469; L=3546; M=3922;
        N=2938; O=984; P=84; Q=683; R=-663; S=229; T=185; U=498; V=1910;
        W=15.13; X=12.65; Y=10.26; Z=7.86; AA=90.42; AB=3437.57;
        AC=1674; AD=24.38; AE=17643; AF=2.31; AG=185; AH=0.45; AI=10.99; AJ=16748936
    }
    6 = @{
        D=5617; E=835; F=948; G=1113; H=902; I=627; K=8817; L=3863; M=4954; N=3814;
        P=90; Q=385; R=-740; S=545; T=161; U=223; V=1840;
        W=14.86; X=16.06; Y=18.57; Z=11.08; AA=77.97; AB=4191.73;
        AC=3448; AD=8.57; AE=19657; AF=1.5; AG=195; AH=0.66; AI=6.54; AJ=17938966
    }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($c in $cols.Keys) {
        $ws.Range("$c$r").Value = $cols[$c]
    }
}

# Estimate-year rows lose all of their figures (only id / label columns A-C survive)
foreach ($r in 7,8,9) {
    $ws.Range("D${r}:AJ${r}").ClearContents()
}
